$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New simulated stock price values for B2:B91 (S0 re-simulated from 40 to 49)
$values = @(
    49.46212935676843,
    49.33602558098684,
    49.94662382203796,
    51.40819847409598,
    51.18429634055894,
    50.96138519705672,
    52.50837614737682,
    53.27866919244317,
    52.81123799998915,
    53.35873116925643,
    52.89665181037189,
    52.43628280325093,
    52.67977582013763,
    50.81208479764181,
    49.18538495699213,
    48.66842598591379,
    47.7486104684644,
    48.03591092007464,
    47.22148718216683,
    45.9805824775968,
    47.27480728948368,
    47.07635979114276,
    47.13931260866798,
    45.88977375508578,
    45.42282011073564,
    45.52087704596488,
    44.54409134900789,
    44.86418133387139,
    44.36047727624057,
    44.11926876524473,
    43.62304766969089,
    45.17982271105376,
    45.17101888191685,
    44.27973135855528,
    44.97612850137143,
    43.95297390773155,
    44.1294619620521,
    42.52761116391168,
    41.47593481417574,
    41.63303219467908,
    42.22066505917076,
    42.36017039404746,
    42.27023562123649,
    42.0329172932786,
    40.87717381264002,
    40.32728252110365,
    39.98015330620286,
    40.78935287917128,
    41.05756294875189,
    39.71452230227865,
    39.96090998109428,
    39.67354319905203,
    39.17160478707117,
    39.62942526457849,
    40.41157107632936,
    41.13156737429249,
    40.48681483589817,
    40.25334030966025,
    40.50856659103023,
    41.2647958244423,
    40.89526907510976,
    40.75448738803336,
    39.91365038982944,
    39.02382639236231,
    39.63005836912316,
    40.66137399995906,
    40.6085213188974,
    41.38850192796388,
    41.67483969780175,
    41.17231882461203,
    41.45697261619378,
    42.68221985103808,
    42.65589919147479,
    43.93866829121426,
    41.81881771478086,
    42.47599584642117,
    42.54848887702006,
    42.31128456960304,
    42.38727221066272,
    40.82712970612457,
    40.66043832369331,
    40.93823937176331,
    42.10027779099923,
    41.69245366440129,
    41.06274426477556,
    40.67766297854446,
    41.38997853644605,
    41.65043412817278,
    41.2380122163849,
    41.64246721067363
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $values[$i]
}
